$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5283
$ws1.Range("F4").Value = 10738
$ws1.Range("F7").Value = 150
$ws1.Range("F8").Value = 185
$ws1.Range("F9").Value = 895

# Sheet "全部类型" (all types) - same events repeated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5283
$ws4.Range("F7").Value = 10738
$ws4.Range("F10").Value = 150
$ws4.Range("F13").Value = 185
$ws4.Range("F14").Value = 895
